# Auto-generated edit script applying numeric updates from the commit diff
# to the Jenova_Profits workbook (sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 111117880
$ws.Range("N76").Value = -166673960
$ws.Range("J76").Value = 166673330
$ws.Range("L76").Value = 166673330
$ws.Range("L79").Value = 166673330
$ws.Range("H79").Value = 111117880
$ws.Range("J79").Value = 166673330
$ws.Range("N79").Value = -166675514
$ws.Range("H112").Value = 3346.9285
$ws.Range("J112").Value = 3530.9614
$ws.Range("N112").Value = -12808.8842
$ws.Range("L112").Value = 10592.8842
$ws.Range("L137").Value = 2738972.4
$ws.Range("N137").Value = -2744072.4
$ws.Range("K137").Value = 1673443.5
$ws.Range("I137").Value = 557814.5
$ws.Range("J137").Value = 912990.8
$ws.Range("M137").Value = -1670893.5
$ws.Range("H137").Value = 692536.5600000001
$ws.Range("J138").Value = 5642.922
$ws.Range("H138").Value = 4380.187
$ws.Range("N138").Value = -27208.766
$ws.Range("L138").Value = 16928.766

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("J2").Value = 844.4
$ws.Range("K2").Value = 1285.381
$ws.Range("N2").Value = -1070.4
$ws.Range("L2").Value = 844.4
$ws.Range("M2").Value = -1172.381
$ws.Range("I2").Value = 1285.381
$ws.Range("H2").Value = 1200.5769
$ws.Range("N86").Value = -52372
$ws.Range("H86").Value = 50000
$ws.Range("J86").Value = 50000
$ws.Range("L86").Value = 50000
$ws.Range("J89").Value = 50000
$ws.Range("L89").Value = 150000
$ws.Range("N89").Value = -161856
$ws.Range("H89").Value = 50000
$ws.Range("M97").Value = -262.2353000000001
$ws.Range("H97").Value = 828.9545000000001
$ws.Range("K97").Value = 758.2353000000001
$ws.Range("I97").Value = 758.2353000000001
$ws.Range("H102").Value = 1779.8572
$ws.Range("K102").Value = 1751.6666
$ws.Range("M102").Value = -129.6666
$ws.Range("I102").Value = 1751.6666
$ws.Range("L110").Value = 625.6
$ws.Range("I110").Value = 129465.46
$ws.Range("J110").Value = 625.6
$ws.Range("N110").Value = -4715.6
$ws.Range("M110").Value = -127420.46
$ws.Range("K110").Value = 129465.46
$ws.Range("H110").Value = 114824.57
$ws.Range("K116").Value = 1285.381
$ws.Range("H116").Value = 1200.5769
$ws.Range("M116").Value = 1008.619
$ws.Range("L116").Value = 844.4
$ws.Range("N116").Value = -5432.4
$ws.Range("J116").Value = 844.4
$ws.Range("I116").Value = 1285.381
$ws.Range("J128").Value = 54966.668
$ws.Range("H128").Value = 54966.668
$ws.Range("N128").Value = -64926.668
$ws.Range("L128").Value = 54966.668

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N3").Value = -1072.4
$ws.Range("H3").Value = 1200.5769
$ws.Range("L3").Value = 844.4
$ws.Range("J3").Value = 844.4
$ws.Range("M3").Value = -1171.381
$ws.Range("K3").Value = 1285.381
$ws.Range("I3").Value = 1285.381
$ws.Range("N28").ClearContents()
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("K86").Value = 1134893.6
$ws.Range("I86").Value = 1134893.6
$ws.Range("M86").Value = -1133770.6
$ws.Range("H86").Value = 851445.7
$ws.Range("K89").Value = 5674468
$ws.Range("H89").Value = 851445.7
$ws.Range("M89").Value = -5668852
$ws.Range("I89").Value = 1134893.6
$ws.Range("H94").Value = 1288.8889
$ws.Range("J94").Value = 1320
$ws.Range("L94").Value = 1320
$ws.Range("N94").Value = -2222
$ws.Range("H99").Value = 1151.4546
$ws.Range("I99").Value = 1179.5555
$ws.Range("M99").Value = 318.4445000000001
$ws.Range("K99").Value = 1179.5555
$ws.Range("K107").Value = 1367.9
$ws.Range("J107").Value = 2002565.2
$ws.Range("H107").Value = 668433.7
$ws.Range("M107").Value = 552.0999999999999
$ws.Range("N107").Value = -2006405.2
$ws.Range("L107").Value = 2002565.2
$ws.Range("I107").Value = 1367.9
$ws.Range("I134").Value = 1573.0667
$ws.Range("M134").Value = -2184.2001
$ws.Range("H134").Value = 22025.256
$ws.Range("K134").Value = 4719.2001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M16").Value = -606.625
$ws.Range("I16").Value = 893.625
$ws.Range("K16").Value = 893.625
$ws.Range("H16").Value = 1102
$ws.Range("J58").Value = 6341.8335
$ws.Range("N58").Value = -6747.8335
$ws.Range("L58").Value = 6341.8335
$ws.Range("H58").Value = 268033.25
$ws.Range("K58").Value = 388813.94
$ws.Range("M58").Value = -388610.94
$ws.Range("I58").Value = 388813.94
$ws.Range("H106").Value = 99500
$ws.Range("L106").Value = 99500
$ws.Range("N106").Value = -102024
$ws.Range("J106").Value = 99500
$ws.Range("K107").Value = 500.16666
$ws.Range("H107").Value = 571.4286
$ws.Range("M107").Value = 1419.83334
$ws.Range("I107").Value = 500.16666
$ws.Range("K113").Value = 893.625
$ws.Range("H113").Value = 1102
$ws.Range("M113").Value = 1276.375
$ws.Range("I113").Value = 893.625
$ws.Range("M122").Value = -4555.6666
$ws.Range("K122").Value = 7005.6666
$ws.Range("I122").Value = 2335.2222
$ws.Range("H122").Value = 4291
$ws.Range("L132").Value = 16849.9995
$ws.Range("H132").Value = 3893.625
$ws.Range("N132").Value = -21909.9995
$ws.Range("K132").Value = 8579.400000000001
$ws.Range("J132").Value = 5616.6665
$ws.Range("I132").Value = 2859.8
$ws.Range("M132").Value = -6049.400000000001
$ws.Range("M136").Value = -1163891.82
$ws.Range("H136").Value = 268033.25
$ws.Range("K136").Value = 1166441.82
$ws.Range("N136").Value = -24125.5005
$ws.Range("L136").Value = 19025.5005
$ws.Range("I136").Value = 388813.94
$ws.Range("J136").Value = 6341.8335

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 350
$ws.Range("J60").Value = 550
$ws.Range("N60").Value = -2152
$ws.Range("L60").Value = 1650
$ws.Range("J122").Value = 111863.664
$ws.Range("M122").Value = -4318
$ws.Range("K122").Value = 6768
$ws.Range("N122").Value = -1011672.976
$ws.Range("I122").Value = 752
$ws.Range("H122").Value = 67419
$ws.Range("L122").Value = 1006772.976

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 44285
$ws.Range("J94").Value = 44285
$ws.Range("L94").Value = 44285
$ws.Range("N94").Value = -45637
$ws.Range("I126").Value = 2557.2856
$ws.Range("J126").Value = 4428.143
$ws.Range("L126").Value = 13284.429
$ws.Range("M126").Value = -5201.8568
$ws.Range("K126").Value = 7671.8568
$ws.Range("H126").Value = 3492.7144
$ws.Range("N126").Value = -18224.429

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("L61").Value = 5442
$ws.Range("N61").Value = -5846
$ws.Range("J61").Value = 5442
$ws.Range("H61").Value = 5060.2856
$ws.Range("L113").Value = 5442
$ws.Range("H113").Value = 5060.2856
$ws.Range("J113").Value = 5442
$ws.Range("N113").Value = -9782

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 7710.3335
$ws.Range("M41").ClearContents()
$ws.Range("K41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("N74").Value = -9510.25
$ws.Range("H74").Value = 6825.273
$ws.Range("J74").Value = 7638.25
$ws.Range("K74").Value = 4657.3335
$ws.Range("I74").Value = 4657.3335
$ws.Range("L74").Value = 7638.25
$ws.Range("M74").Value = -3721.3335
$ws.Range("I77").Value = 4657.3335
$ws.Range("M77").Value = -9292.000499999998
$ws.Range("L77").Value = 22914.75
$ws.Range("N77").Value = -32274.75
$ws.Range("J77").Value = 7638.25
$ws.Range("H77").Value = 6825.273
$ws.Range("K77").Value = 13972.0005
$ws.Range("I96").Value = 334731.66
$ws.Range("M96").Value = -333358.66
$ws.Range("L96").Value = 1999.5
$ws.Range("H96").Value = 201638.8
$ws.Range("J96").Value = 1999.5
$ws.Range("K96").Value = 334731.66
$ws.Range("N96").Value = -4745.5
$ws.Range("L113").Value = 1109.1429
$ws.Range("K113").Value = 2075.0001
$ws.Range("H113").Value = 466.3
$ws.Range("J113").Value = 369.7143
$ws.Range("M113").Value = 94.9998999999998
$ws.Range("N113").Value = -5449.1429
$ws.Range("I113").Value = 691.6667
$ws.Range("I126").Value = 763.7778
$ws.Range("M126").Value = 178.6666
$ws.Range("K126").Value = 2291.3334
$ws.Range("H126").Value = 737.4

